$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (risk #4) content was replaced: the old "Sitting on chair too long" /
# ergonomics risk was swapped out for a new "Traffic monitored" / HTTPS risk.
$ws.Range("B6").Value = "Traffic monitored"
$ws.Range("C6").Value = "Any input shown in plain text"
$ws.Range("D6").Value = "Use HTTPS only, however no sensitve information is being entered at this stage"
$ws.Range("E6").Value = "use of secure protocols to avoid data being stolen"
$ws.Range("F6").Value = "Likely"
$ws.Range("H6").Value = 3

# Selection moved
$ws.Range("F17").Select()
